$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B1: "select" -> "delete"
$ws.Range("B1").Value = "delete"

# D4: "include" -> empty
$ws.Range("D4").Value = ""

# H6:H15 - rewrite SELECT (...) statements into DELETE statements
$ws.Range("H6").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-1 IO'  AND  IsConnected = '0';"
$ws.Range("H7").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-2 IO'  AND  IsConnected = '0';"
$ws.Range("H8").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-3 IO'  AND  IsConnected = '0';"
$ws.Range("H9").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-4 IO'  AND  IsConnected = '0';"
$ws.Range("H10").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-5 IO'  AND  IsConnected = '0';"
$ws.Range("H11").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-6 IO'  AND  IsConnected = '0';"
$ws.Range("H12").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-7 IO'  AND  IsConnected = '0';"
$ws.Range("H13").Value = "DELETE FROM IOChannels WHERE Name = 'Charger 7-8 IO'  AND  IsConnected = '0';"
$ws.Range("H14").Value = "DELETE FROM IOChannels WHERE Name = 'Plant 7 IO'  AND  IsConnected = '0';"
$ws.Range("H15").Value = "DELETE FROM IOChannels WHERE Name = 'Charget 7-9 IO'  AND  IsConnected = '0';"

$wb.Save()
